$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: a handful of "Price" (column D) values are digit-and-dot strings that
# Excel would otherwise auto-convert to a number on assignment (e.g. "571.76").
# The source workbook stores every Price/Volume cell as text, so for those cells
# we briefly force a Text number format, assign the string, then clear the format
# again so no stray style index is left behind on the cell.

$ws.Range("D2").Value = '63.160.04'
$ws.Range("E2").Value = '  +0.85%  '
$ws.Range("D3").Value = '2.448.31'
$ws.Range("E3").Value = '  +0.44%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '571.76'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.86%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.47'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.73%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.536'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.73%  '
$ws.Range("D9").Value = '2.444.79'
$ws.Range("E9").Value = '  +0.24%  '
$ws.Range("E10").Value = '  +0.27%  '
$ws.Range("E11").Value = '  +1.28%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.28'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.52%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.355'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.02%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.95'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.34%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000178'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.03%  '
$ws.Range("D17").Value = '62.856.21'
$ws.Range("E17").Value = '  +0.71%  '
$ws.Range("D18").Value = '2.439.73'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.27'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.39%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.33'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +5.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '327.89'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.20'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.79%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.06'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +11.96%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.01'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.72%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.70'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -2.41%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '611.27'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +5.08%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.90'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +3.69%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000102'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.54%  '
$ws.Range("D29").Value = '2.577.79'
$ws.Range("E29").Value = '  +0.78%  '
$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.32%  '
$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.49'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +3.70%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.22'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -2.41%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.141'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -3.10%  '
$ws.Range("E34").Value = '  +1.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.18'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +6.81%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.52'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.16%  '
$ws.Range("E37").Value = '  +0.12%  '
$ws.Range("E38").Value = '  -0.73%  '
$ws.Range("B39").Value = 'EthereumClassic'
$ws.Range("C39").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.76'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.02%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.40'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.29%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '147.36'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.54%  '
$ws.Range("E42").Value = '  -1.98%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.60'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +6.09%  '
$ws.Range("E44").Value = '  -0.15%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '41.98'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.89%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '148.57'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.75'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.40%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '21.13'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +2.78%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.600'
$ws.Range("D50").ClearFormats()
$ws.Range("E51").Value = '  +0.19%  '
